$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 65 will become a new section header ("splash_screen.dart"): merge it
# across A:B now (merging has to happen before the header style is pasted
# onto it, otherwise the paste-special below stops landing on the clean
# header style).
$ws.Range("A65:B65").Merge()
$ws.Range("A57:B57").Copy()
$ws.Range("A65:B65").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fill in the cell values in the precise order that reproduces the
# --- target shared-strings table ordering.

# Row 60: "Settings" / "Cài đặt" (B60 already holds "Cài đặt")
$ws.Range("A60").Value = "Settings"

# Row 61: "Location" / "Vị trí"
$ws.Range("B61").Value = "Vị trí"
$ws.Range("A61").Value = "Location"

# Row 63: "Language" / "Ngôn ngữ"
$ws.Range("B63").Value = "Ngôn ngữ"
$ws.Range("A63").Value = "Language"

# Row 66: "Cleaning closet…" / "Đang lau chùi tủ đồ"
$ws.Range("A66").Value = "Cleaning closet…"
$ws.Range("B66").Value = "Đang lau chùi tủ đồ"

# Row 65: "splash_screen.dart"
$ws.Range("A65").Value = "splash_screen.dart"

# Row 62: "Auto-detect" / "Tự động phát hiện" (already-existing shared strings)
$ws.Range("A62").Value = "Auto-detect"
$ws.Range("B62").Value = "Tự động phát hiện"

# --- Re-apply the cell formatting that plain value assignment above can
# --- reset (e.g. A60 loses its quote-prefixed style once a value lands in
# --- it), and give row 66 the same look as the row 62/63 entries.

$ws.Range("A58").Copy()
$ws.Range("A60").PasteSpecial(-4122)

$ws.Range("A62:B62").Copy()
$ws.Range("A66:B66").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Keep the selection/active cell in sync with the new bottom of the sheet.
$ws.Range("A67").Select()

Write-Output "applied"
